$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "City"
$ws.Range("E1").Value = "Email"
$ws.Range("F1").Value = "Contact Person"
$ws.Range("G1").Value = "Contact Mobile"

$ws.Range("H15").Select()
